$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted for Espárragos (row 39), pushing all
# subsequent rows (formerly 39-54) down by one (now 40-55).
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly entry.
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 44523
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 300000000
$ws.Range("G39").Value = "Espárragos"
$ws.Range("H39").Value = "Verde"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 1100
$ws.Range("L39").Value = 1100
$ws.Range("M39").Value = 1100
$ws.Range("N39").Value = "`$/kilo"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 1100
$ws.Range("Q39").Value = 1
$ws.Range("R39").Value = "Hortaliza"
